$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Era C, Expense)
$ws.Range("C2").Value = -245824.33
$ws.Range("D2").Value = -245225.59
$ws.Range("H2").Value = -491049.92

# Row 3 (Era C, Income)
$ws.Range("C3").Value = 591845.34
$ws.Range("D3").Value = 585293.16
$ws.Range("H3").Value = 1177138.5

# Row 5 (Era B, Expense)
$ws.Range("E5").Value = -237152.43
$ws.Range("H5").Value = -237152.43

# Row 6 (Era B, Income)
$ws.Range("E6").Value = 643798.41
$ws.Range("H6").Value = 643798.41

# Row 8 (Era A, Expense)
$ws.Range("F8").Value = -265602.84
$ws.Range("G8").Value = -221532.92
$ws.Range("H8").Value = -487135.76

# Row 9 (Era A, Income)
$ws.Range("F9").Value = 537409.59
$ws.Range("G9").Value = 647993.1
$ws.Range("H9").Value = 1185402.69

# Row 11 (Total)
$ws.Range("C11").Value = 346021.01
$ws.Range("D11").Value = 340067.57
$ws.Range("E11").Value = 406645.98
$ws.Range("F11").Value = 271806.75
$ws.Range("G11").Value = 426460.18
$ws.Range("H11").Value = 1791001.49
